{"js": "// Move the \"2022\" year token from mid-sentence to the start of the\n// observation-dates sentence (appears 4 times in the document body).\nconst oldText =\n  \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03bc\u03c0\u03cc\u03c4\u03b5\u03c2 2022: 14-23 \u039c\u03b1\u0390\u03bf\u03c5, 13-22 \u0399\u03bf\u03c5\u03bd\u03af\u03bf\u03c5, 12-21 \u0399\u03bf\u03c5\u03bb\u03af\u03bf\u03c5\";\nconst newText =\n  \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03bc\u03c0\u03cc\u03c4\u03b5\u03c2: 14-23 \u039c\u03b1\u0390\u03bf\u03c5, 13-22 \u0399\u03bf\u03c5\u03bd\u03af\u03bf\u03c5, 12-21 \u0399\u03bf\u03c5\u03bb\u03af\u03bf\u03c5\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nresults.items.forEach((range) => {\n  range.insertText(newText, \"Replace\");\n});\n\nawait context.sync();\n", "ps1": "# Move the \"2022\" year token from mid-sentence to the start of the\n# observation-dates sentence (appears 4 times in the document body).\n$d = $word.ActiveDocument\n\n$oldText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03bc\u03c0\u03cc\u03c4\u03b5\u03c2 2022: 14-23 \u039c\u03b1\u0390\u03bf\u03c5, 13-22 \u0399\u03bf\u03c5\u03bd\u03af\u03bf\u03c5, 12-21 \u0399\u03bf\u03c5\u03bb\u03af\u03bf\u03c5\"\n$newText = \"2022 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 \u0391\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc\u03c2 \u03bc\u03c0\u03cc\u03c4\u03b5\u03c2: 14-23 \u039c\u03b1\u0390\u03bf\u03c5, 13-22 \u0399\u03bf\u03c5\u03bd\u03af\u03bf\u03c5, 12-21 \u0399\u03bf\u03c5\u03bb\u03af\u03bf\u03c5\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n$find.Execute(\n    $oldText,   # FindText\n    $true,      # MatchCase\n    $false,     # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue) - scan the whole story\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n)\n"}
